# Weekly update: a new price record is added for "Poroto verde" (Hortaliza)
# at the Mercado Mayorista Lo Valledor de Santiago market. The new record is
# inserted as row 1232, pushing all the existing rows from 1232..1337 down by
# one (to 1233..1338), exactly like a normal spreadsheet row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 1232.
$ws.Rows.Item(1232).Insert()

# Populate the newly inserted row with the latest weekly price data.
$ws.Range("A1232").Value = 6
$ws.Range("B1232").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1232").Value = "Metropolitana"
$ws.Range("D1232").Value = 45132
$ws.Range("E1232").Value = 13
$ws.Range("F1232").Value = 100112031
$ws.Range("G1232").Value = "Poroto verde"
$ws.Range("H1232").Value = "Magnum"
$ws.Range("I1232").Value = "Primera"
$ws.Range("J1232").Value = 460
$ws.Range("K1232").Value = 17000
$ws.Range("L1232").Value = 18000
$ws.Range("M1232").Value = 17500
$ws.Range("N1232").Value = "`$/malla 25 kilos"
$ws.Range("O1232").Value = "Perú"
$ws.Range("P1232").Value = 700
$ws.Range("Q1232").Value = 25
$ws.Range("R1232").Value = "Hortaliza"
